$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.140.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.750.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.51"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5287"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2818"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06189"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.746.83"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07188"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.48"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6463"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.635"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.63"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.049.14"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.78"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006743"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.971.57"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.334"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.747"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.237"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.519"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.32"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.814"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08306"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.806"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.648"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04624"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.648"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.018"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6356"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.06%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01625"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.990"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9995"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.73"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3938"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7526"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.064"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1155"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.377"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05353"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.67"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.07"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3490"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.602"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.70%  "
